$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("A2").Value = "ECs"
    $ws.Range("B2").Value = "Rtn4"
    $ws.Range("C2").Value = "Lingo1"
    $ws.Range("D2").Value = "FAPs"
    $ws.Range("E2").Value = 2
    $ws.Range("F2").Value = 1
    $ws.Range("G2").Value = 78.25435300000001
    $ws.Range("H2").Value = 156.508706
    $ws.Range("I2").Value = 0.1800096451968904
    $ws.Range("J2").Value = 0.1347002616492669
    $ws.Range("K2").Value = 1
    $ws.Range("L2").Value = 0.3333333333333333
    $ws.Range("M2").Value = 0.096568
    $ws.Range("N2").Value = 0.289704
    $ws.Range("O2").Value = 0.3543662778110081
    $ws.Range("P2").Value = 0.4515435187176778
    $ws.Range("Q2").Value = 7.556866360504001
    $ws.Range("R2").Value = 45.34119816302401
    $ws.Range("S2").Value = 0.06378934793850224
    $ws.Range("T2").Value = 0.06082303011730185
    # Row 3
    $ws.Range("A3").Value = "ECs"
    $ws.Range("B3").Value = "Rtn4"
    $ws.Range("C3").Value = "Lingo1"
    $ws.Range("D3").Value = "MuSCs"
    $ws.Range("E3").Value = 2
    $ws.Range("F3").Value = 1
    $ws.Range("G3").Value = 78.25435300000001
    $ws.Range("H3").Value = 156.508706
    $ws.Range("I3").Value = 0.1800096451968904
    $ws.Range("J3").Value = 0.1347002616492669
    $ws.Range("K3").Value = 1
    $ws.Range("L3").Value = 0.5
    $ws.Range("M3").Value = 0.175941
    $ws.Range("N3").Value = 0.351882
    $ws.Range("O3").Value = 0.6456337221889918
    $ws.Range("P3").Value = 0.5484564812823222
    $ws.Range("Q3").Value = 13.768149121173
    $ws.Range("R3").Value = 55.072596484692
    $ws.Range("S3").Value = 0.1162202972583881
    $ws.Range("T3").Value = 0.07387723153196506
    # Row 4
    $ws.Range("A4").Value = "FAPs"
    $ws.Range("B4").Value = "Rtn4"
    $ws.Range("C4").Value = "Lingo1"
    $ws.Range("D4").Value = "FAPs"
    $ws.Range("E4").Value = 3
    $ws.Range("F4").Value = 1
    $ws.Range("G4").Value = 80.25908033333333
    $ws.Range("H4").Value = 240.777241
    $ws.Range("I4").Value = 0.1846211491216604
    $ws.Range("J4").Value = 0.2072265383236163
    $ws.Range("K4").Value = 1
    $ws.Range("L4").Value = 0.3333333333333333
    $ws.Range("M4").Value = 0.096568
    $ws.Range("N4").Value = 0.289704
    $ws.Range("O4").Value = 0.3543662778110081
    $ws.Range("P4").Value = 0.4515435187176778
    $ws.Range("Q4").Value = 7.750458869629333
    $ws.Range("R4").Value = 69.754129826664
    $ws.Range("S4").Value = 0.06542350941943385
    $ws.Range("T4").Value = 0.09357180028632939
    # Row 5
    $ws.Range("A5").Value = "FAPs"
    $ws.Range("B5").Value = "Rtn4"
    $ws.Range("C5").Value = "Lingo1"
    $ws.Range("D5").Value = "MuSCs"
    $ws.Range("E5").Value = 3
    $ws.Range("F5").Value = 1
    $ws.Range("G5").Value = 80.25908033333333
    $ws.Range("H5").Value = 240.777241
    $ws.Range("I5").Value = 0.1846211491216604
    $ws.Range("J5").Value = 0.2072265383236163
    $ws.Range("K5").Value = 1
    $ws.Range("L5").Value = 0.5
    $ws.Range("M5").Value = 0.175941
    $ws.Range("N5").Value = 0.351882
    $ws.Range("O5").Value = 0.6456337221889918
    $ws.Range("P5").Value = 0.5484564812823222
    $ws.Range("Q5").Value = 14.120862852927
    $ws.Range("R5").Value = 84.725177117562
    $ws.Range("S5").Value = 0.1191976397022265
    $ws.Range("T5").Value = 0.1136547380372869
    # Row 6
    $ws.Range("A6").Value = "Inflammatory-Mac"
    $ws.Range("B6").Value = "Rtn4"
    $ws.Range("C6").Value = "Lingo1"
    $ws.Range("D6").Value = "FAPs"
    $ws.Range("E6").Value = 3
    $ws.Range("F6").Value = 1
    $ws.Range("G6").Value = 73.92583733333333
    $ws.Range("H6").Value = 221.777512
    $ws.Range("I6").Value = 0.1700526966117318
    $ws.Range("J6").Value = 0.1908742948416136
    $ws.Range("K6").Value = 1
    $ws.Range("L6").Value = 0.3333333333333333
    $ws.Range("M6").Value = 0.096568
    $ws.Range("N6").Value = 0.289704
    $ws.Range("O6").Value = 0.3543662778110081
    $ws.Range("P6").Value = 0.4515435187176778
    $ws.Range("Q6").Value = 7.138870259605334
    $ws.Range("R6").Value = 64.24983233644801
    $ws.Range("S6").Value = 0.06026094113002403
    $ws.Range("T6").Value = 0.08618805072553772
    # Row 7
    $ws.Range("A7").Value = "Inflammatory-Mac"
    $ws.Range("B7").Value = "Rtn4"
    $ws.Range("C7").Value = "Lingo1"
    $ws.Range("D7").Value = "MuSCs"
    $ws.Range("E7").Value = 3
    $ws.Range("F7").Value = 1
    $ws.Range("G7").Value = 73.92583733333333
    $ws.Range("H7").Value = 221.777512
    $ws.Range("I7").Value = 0.1700526966117318
    $ws.Range("J7").Value = 0.1908742948416136
    $ws.Range("K7").Value = 1
    $ws.Range("L7").Value = 0.5
    $ws.Range("M7").Value = 0.175941
    $ws.Range("N7").Value = 0.351882
    $ws.Range("O7").Value = 0.6456337221889918
    $ws.Range("P7").Value = 0.5484564812823222
    $ws.Range("Q7").Value = 13.006585746264
    $ws.Range("R7").Value = 78.03951447758399
    $ws.Range("S7").Value = 0.1097917554817078
    $ws.Range("T7").Value = 0.1046862441160759
    # Row 8
    $ws.Range("A8").Value = "MuSCs"
    $ws.Range("B8").Value = "Rtn4"
    $ws.Range("C8").Value = "Lingo1"
    $ws.Range("D8").Value = "FAPs"
    $ws.Range("E8").Value = 2
    $ws.Range("F8").Value = 1
    $ws.Range("G8").Value = 64.0114765
    $ws.Range("H8").Value = 128.022953
    $ws.Range("I8").Value = 0.1472465458029419
    $ws.Range("J8").Value = 0.1101838083448968
    $ws.Range("K8").Value = 1
    $ws.Range("L8").Value = 0.3333333333333333
    $ws.Range("M8").Value = 0.096568
    $ws.Range("N8").Value = 0.289704
    $ws.Range("O8").Value = 0.3543662778110081
    $ws.Range("P8").Value = 0.4515435187176778
    $ws.Range("Q8").Value = 6.181460262652
    $ws.Range("R8").Value = 37.088761575912
    $ws.Range("S8").Value = 0.05217921035671665
    $ws.Range("T8").Value = 0.04975278452576892
    # Row 9
    $ws.Range("A9").Value = "MuSCs"
    $ws.Range("B9").Value = "Rtn4"
    $ws.Range("C9").Value = "Lingo1"
    $ws.Range("D9").Value = "MuSCs"
    $ws.Range("E9").Value = 2
    $ws.Range("F9").Value = 1
    $ws.Range("G9").Value = 64.0114765
    $ws.Range("H9").Value = 128.022953
    $ws.Range("I9").Value = 0.1472465458029419
    $ws.Range("J9").Value = 0.1101838083448968
    $ws.Range("K9").Value = 1
    $ws.Range("L9").Value = 0.5
    $ws.Range("M9").Value = 0.175941
    $ws.Range("N9").Value = 0.351882
    $ws.Range("O9").Value = 0.6456337221889918
    $ws.Range("P9").Value = 0.5484564812823222
    $ws.Range("Q9").Value = 11.2622431868865
    $ws.Range("R9").Value = 45.048972747546
    $ws.Range("S9").Value = 0.09506733544622527
    $ws.Range("T9").Value = 0.06043102381912786
    # Row 10
    $ws.Range("A10").Value = "Neutrophils"
    $ws.Range("B10").Value = "Rtn4"
    $ws.Range("C10").Value = "Lingo1"
    $ws.Range("D10").Value = "FAPs"
    $ws.Range("E10").Value = 3
    $ws.Range("F10").Value = 1
    $ws.Range("G10").Value = 56.48610166666666
    $ws.Range("H10").Value = 169.458305
    $ws.Range("I10").Value = 0.1299358148111217
    $ws.Range("J10").Value = 0.1458454203955994
    $ws.Range("K10").Value = 1
    $ws.Range("L10").Value = 0.3333333333333333
    $ws.Range("M10").Value = 0.096568
    $ws.Range("N10").Value = 0.289704
    $ws.Range("O10").Value = 0.3543662778110081
    $ws.Range("P10").Value = 0.4515435187176778
    $ws.Range("Q10").Value = 5.454749865746666
    $ws.Range("R10").Value = 49.09274879172
    $ws.Range("S10").Value = 0.04604487104895764
    $ws.Range("T10").Value = 0.06585555431428791
    # Row 11
    $ws.Range("A11").Value = "Neutrophils"
    $ws.Range("B11").Value = "Rtn4"
    $ws.Range("C11").Value = "Lingo1"
    $ws.Range("D11").Value = "MuSCs"
    $ws.Range("E11").Value = 3
    $ws.Range("F11").Value = 1
    $ws.Range("G11").Value = 56.48610166666666
    $ws.Range("H11").Value = 169.458305
    $ws.Range("I11").Value = 0.1299358148111217
    $ws.Range("J11").Value = 0.1458454203955994
    $ws.Range("K11").Value = 1
    $ws.Range("L11").Value = 0.5
    $ws.Range("M11").Value = 0.175941
    $ws.Range("N11").Value = 0.351882
    $ws.Range("O11").Value = 0.6456337221889918
    $ws.Range("P11").Value = 0.5484564812823222
    $ws.Range("Q11").Value = 9.938221213334998
    $ws.Range("R11").Value = 59.62932728001
    $ws.Range("S11").Value = 0.08389094376216401
    $ws.Range("T11").Value = 0.07998986608131146
    # Row 12
    $ws.Range("A12").Value = "Resolving-Mac"
    $ws.Range("B12").Value = "Rtn4"
    $ws.Range("C12").Value = "Lingo1"
    $ws.Range("D12").Value = "FAPs"
    $ws.Range("E12").Value = 3
    $ws.Range("F12").Value = 1
    $ws.Range("G12").Value = 81.78626233333334
    $ws.Range("H12").Value = 245.358787
    $ws.Range("I12").Value = 0.1881341484556537
    $ws.Range("J12").Value = 0.2111696764450071
    $ws.Range("K12").Value = 1
    $ws.Range("L12").Value = 0.3333333333333333
    $ws.Range("M12").Value = 0.096568
    $ws.Range("N12").Value = 0.289704
    $ws.Range("O12").Value = 0.3543662778110081
    $ws.Range("P12").Value = 0.4515435187176778
    $ws.Range("Q12").Value = 7.897935781005334
    $ws.Range("R12").Value = 71.081422029048
    $ws.Range("S12").Value = 0.06666839791737364
    $ws.Range("T12").Value = 0.09535229874845204
    # Row 13
    $ws.Range("A13").Value = "Resolving-Mac"
    $ws.Range("B13").Value = "Rtn4"
    $ws.Range("C13").Value = "Lingo1"
    $ws.Range("D13").Value = "MuSCs"
    $ws.Range("E13").Value = 3
    $ws.Range("F13").Value = 1
    $ws.Range("G13").Value = 81.78626233333334
    $ws.Range("H13").Value = 245.358787
    $ws.Range("I13").Value = 0.1881341484556537
    $ws.Range("J13").Value = 0.2111696764450071
    $ws.Range("K13").Value = 1
    $ws.Range("L13").Value = 0.5
    $ws.Range("M13").Value = 0.175941
    $ws.Range("N13").Value = 0.351882
    $ws.Range("O13").Value = 0.6456337221889918
    $ws.Range("P13").Value = 0.5484564812823222
    $ws.Range("Q13").Value = 14.389556781189
    $ws.Range("R13").Value = 86.337340687134
    $ws.Range("S13").Value = 0.1214657505382801
    $ws.Range("T13").Value = 0.1158173776965551
